$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{ Row = 17; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 21; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 23; I = '%'; J = 'Uninterpretable' }
    @{ Row = 32; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 35; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 38; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 40; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 62; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 86; I = 'qy'; J = 'Yes-No-Question' }
    @{ Row = 127; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 162; I = 'ba'; J = 'Appreciation' }
    @{ Row = 170; I = 'qy'; J = 'Yes-No-Question' }
    @{ Row = 185; I = 'ba'; J = 'Appreciation' }
    @{ Row = 206; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 213; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 215; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 220; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 223; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 224; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 246; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 247; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 268; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 270; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 275; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 278; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 280; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 283; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 311; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 347; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 351; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 353; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 362; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 363; I = '%'; J = 'Uninterpretable' }
    @{ Row = 365; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 372; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 388; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 389; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 400; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 401; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 412; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 421; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 422; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 426; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 432; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 440; I = '%'; J = 'Uninterpretable' }
    @{ Row = 458; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 481; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 488; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 514; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 532; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 536; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 550; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 551; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 553; I = '%'; J = 'Uninterpretable' }
    @{ Row = 557; I = '%'; J = 'Uninterpretable' }
    @{ Row = 560; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 583; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 653; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 656; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 662; I = 'sd'; J = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
